{"js": "// Apply 100 find-and-replace edits to update arithmetic expressions in the table.\n// Pairs are processed in document order so that substring ambiguities\n// (e.g. \"2+16=\" being contained in \"22+16=\") resolve correctly: the\n// longer/earlier-occurring expression is replaced before the shorter one\n// that could otherwise spuriously match inside it.\nconst pairs = [\n  [\"14+11=\", \"71+18=\"],\n  [\"50-32=\", \"22+38=\"],\n  [\"96-77=\", \"38-38=\"],\n  [\"24+74=\", \"13+21=\"],\n  [\"72-44=\", \"5+15=\"],\n  [\"54-30=\", \"21+42=\"],\n  [\"74-25=\", \"19+24=\"],\n  [\"25+18=\", \"43+55=\"],\n  [\"91-11=\", \"50-16=\"],\n  [\"95-51=\", \"61-39=\"],\n  [\"8+38=\", \"50+16=\"],\n  [\"50+49=\", \"14+17=\"],\n  [\"45+22=\", \"23+42=\"],\n  [\"54+40=\", \"12+53=\"],\n  [\"8+24=\", \"10+47=\"],\n  [\"22+16=\", \"39-4=\"],\n  [\"47+43=\", \"88-68=\"],\n  [\"2+17=\", \"38-29=\"],\n  [\"87-33=\", \"52+31=\"],\n  [\"30+0=\", \"16+46=\"],\n  [\"58+2=\", \"14+42=\"],\n  [\"37+43=\", \"11-0=\"],\n  [\"33-19=\", \"78+5=\"],\n  [\"27+56=\", \"36-35=\"],\n  [\"67-49=\", \"63+33=\"],\n  [\"63-1=\", \"88-0=\"],\n  [\"68+14=\", \"52+32=\"],\n  [\"33+39=\", \"7-7=\"],\n  [\"66-21=\", \"2+60=\"],\n  [\"87-69=\", \"29+59=\"],\n  [\"47+44=\", \"97-72=\"],\n  [\"77-28=\", \"53+41=\"],\n  [\"17+66=\", \"94-78=\"],\n  [\"39+5=\", \"72+7=\"],\n  [\"12+35=\", \"9+4=\"],\n  [\"70-34=\", \"81-49=\"],\n  [\"86-6=\", \"57-31=\"],\n  [\"25-11=\", \"95-52=\"],\n  [\"34+18=\", \"14+47=\"],\n  [\"93-90=\", \"1+14=\"],\n  [\"43-14=\", \"77-27=\"],\n  [\"63-15=\", \"18-12=\"],\n  [\"24+40=\", \"96-95=\"],\n  [\"46-3=\", \"2+86=\"],\n  [\"82-61=\", \"68+1=\"],\n  [\"57-32=\", \"44-33=\"],\n  [\"80+1=\", \"20+3=\"],\n  [\"72-7=\", \"21+43=\"],\n  [\"57-1=\", \"43+8=\"],\n  [\"66+24=\", \"96-12=\"],\n  [\"2+16=\", \"72-43=\"],\n  [\"23-8=\", \"21+63=\"],\n  [\"51+19=\", \"59-20=\"],\n  [\"37+20=\", \"72-0=\"],\n  [\"81-48=\", \"31+23=\"],\n  [\"99-19=\", \"77-60=\"],\n  [\"38+18=\", \"49+1=\"],\n  [\"70+9=\", \"70-12=\"],\n  [\"32+6=\", \"72-39=\"],\n  [\"78-62=\", \"34-21=\"],\n  [\"48+13=\", \"59+15=\"],\n  [\"20+26=\", \"36+0=\"],\n  [\"68+19=\", \"12+8=\"],\n  [\"83+15=\", \"7+79=\"],\n  [\"41+44=\", \"50-43=\"],\n  [\"43-27=\", \"49+25=\"],\n  [\"28-22=\", \"81+8=\"],\n  [\"46+5=\", \"22+23=\"],\n  [\"40+53=\", \"59+38=\"],\n  [\"77-51=\", \"20+60=\"],\n  [\"64-53=\", \"6+84=\"],\n  [\"11+81=\", \"56+42=\"],\n  [\"49-24=\", \"49-10=\"],\n  [\"16+52=\", \"55+7=\"],\n  [\"20+75=\", \"49-36=\"],\n  [\"65-10=\", \"50-15=\"],\n  [\"50+15=\", \"44-22=\"],\n  [\"43+47=\", \"34+12=\"],\n  [\"35-3=\", \"86-59=\"],\n  [\"40-31=\", \"72-26=\"],\n  [\"74-18=\", \"20+39=\"],\n  [\"13+20=\", \"78-12=\"],\n  [\"41+5=\", \"18-2=\"],\n  [\"37+25=\", \"93-4=\"],\n  [\"27+9=\", \"28+14=\"],\n  [\"32+14=\", \"67-54=\"],\n  [\"90-58=\", \"59+5=\"],\n  [\"63+31=\", \"74-66=\"],\n  [\"77-64=\", \"70+18=\"],\n  [\"70-17=\", \"87-79=\"],\n  [\"41+52=\", \"51-16=\"],\n  [\"92+6=\", \"16+28=\"],\n  [\"94-43=\", \"96-68=\"],\n  [\"26+19=\", \"84-55=\"],\n  [\"83-58=\", \"5+70=\"],\n  [\"39+50=\", \"19+65=\"],\n  [\"16+6=\", \"1+86=\"],\n  [\"82-28=\", \"21+74=\"],\n  [\"62+31=\", \"1+30=\"],\n  [\"2+80=\", \"76+14=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first (and, by construction, only remaining) match.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply 100 find-and-replace edits to update arithmetic expressions in the table.\n# Pairs are processed in document order so that the one substring ambiguity\n# (\"2+16=\" is contained in \"22+16=\") resolves correctly: \"22+16=\" is replaced\n# before \"2+16=\" is searched for, leaving only the genuine standalone match.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"14+11=\", \"71+18=\"),\n    @(\"50-32=\", \"22+38=\"),\n    @(\"96-77=\", \"38-38=\"),\n    @(\"24+74=\", \"13+21=\"),\n    @(\"72-44=\", \"5+15=\"),\n    @(\"54-30=\", \"21+42=\"),\n    @(\"74-25=\", \"19+24=\"),\n    @(\"25+18=\", \"43+55=\"),\n    @(\"91-11=\", \"50-16=\"),\n    @(\"95-51=\", \"61-39=\"),\n    @(\"8+38=\", \"50+16=\"),\n    @(\"50+49=\", \"14+17=\"),\n    @(\"45+22=\", \"23+42=\"),\n    @(\"54+40=\", \"12+53=\"),\n    @(\"8+24=\", \"10+47=\"),\n    @(\"22+16=\", \"39-4=\"),\n    @(\"47+43=\", \"88-68=\"),\n    @(\"2+17=\", \"38-29=\"),\n    @(\"87-33=\", \"52+31=\"),\n    @(\"30+0=\", \"16+46=\"),\n    @(\"58+2=\", \"14+42=\"),\n    @(\"37+43=\", \"11-0=\"),\n    @(\"33-19=\", \"78+5=\"),\n    @(\"27+56=\", \"36-35=\"),\n    @(\"67-49=\", \"63+33=\"),\n    @(\"63-1=\", \"88-0=\"),\n    @(\"68+14=\", \"52+32=\"),\n    @(\"33+39=\", \"7-7=\"),\n    @(\"66-21=\", \"2+60=\"),\n    @(\"87-69=\", \"29+59=\"),\n    @(\"47+44=\", \"97-72=\"),\n    @(\"77-28=\", \"53+41=\"),\n    @(\"17+66=\", \"94-78=\"),\n    @(\"39+5=\", \"72+7=\"),\n    @(\"12+35=\", \"9+4=\"),\n    @(\"70-34=\", \"81-49=\"),\n    @(\"86-6=\", \"57-31=\"),\n    @(\"25-11=\", \"95-52=\"),\n    @(\"34+18=\", \"14+47=\"),\n    @(\"93-90=\", \"1+14=\"),\n    @(\"43-14=\", \"77-27=\"),\n    @(\"63-15=\", \"18-12=\"),\n    @(\"24+40=\", \"96-95=\"),\n    @(\"46-3=\", \"2+86=\"),\n    @(\"82-61=\", \"68+1=\"),\n    @(\"57-32=\", \"44-33=\"),\n    @(\"80+1=\", \"20+3=\"),\n    @(\"72-7=\", \"21+43=\"),\n    @(\"57-1=\", \"43+8=\"),\n    @(\"66+24=\", \"96-12=\"),\n    @(\"2+16=\", \"72-43=\"),\n    @(\"23-8=\", \"21+63=\"),\n    @(\"51+19=\", \"59-20=\"),\n    @(\"37+20=\", \"72-0=\"),\n    @(\"81-48=\", \"31+23=\"),\n    @(\"99-19=\", \"77-60=\"),\n    @(\"38+18=\", \"49+1=\"),\n    @(\"70+9=\", \"70-12=\"),\n    @(\"32+6=\", \"72-39=\"),\n    @(\"78-62=\", \"34-21=\"),\n    @(\"48+13=\", \"59+15=\"),\n    @(\"20+26=\", \"36+0=\"),\n    @(\"68+19=\", \"12+8=\"),\n    @(\"83+15=\", \"7+79=\"),\n    @(\"41+44=\", \"50-43=\"),\n    @(\"43-27=\", \"49+25=\"),\n    @(\"28-22=\", \"81+8=\"),\n    @(\"46+5=\", \"22+23=\"),\n    @(\"40+53=\", \"59+38=\"),\n    @(\"77-51=\", \"20+60=\"),\n    @(\"64-53=\", \"6+84=\"),\n    @(\"11+81=\", \"56+42=\"),\n    @(\"49-24=\", \"49-10=\"),\n    @(\"16+52=\", \"55+7=\"),\n    @(\"20+75=\", \"49-36=\"),\n    @(\"65-10=\", \"50-15=\"),\n    @(\"50+15=\", \"44-22=\"),\n    @(\"43+47=\", \"34+12=\"),\n    @(\"35-3=\", \"86-59=\"),\n    @(\"40-31=\", \"72-26=\"),\n    @(\"74-18=\", \"20+39=\"),\n    @(\"13+20=\", \"78-12=\"),\n    @(\"41+5=\", \"18-2=\"),\n    @(\"37+25=\", \"93-4=\"),\n    @(\"27+9=\", \"28+14=\"),\n    @(\"32+14=\", \"67-54=\"),\n    @(\"90-58=\", \"59+5=\"),\n    @(\"63+31=\", \"74-66=\"),\n    @(\"77-64=\", \"70+18=\"),\n    @(\"70-17=\", \"87-79=\"),\n    @(\"41+52=\", \"51-16=\"),\n    @(\"92+6=\", \"16+28=\"),\n    @(\"94-43=\", \"96-68=\"),\n    @(\"26+19=\", \"84-55=\"),\n    @(\"83-58=\", \"5+70=\"),\n    @(\"39+50=\", \"19+65=\"),\n    @(\"16+6=\", \"1+86=\"),\n    @(\"82-28=\", \"21+74=\"),\n    @(\"62+31=\", \"1+30=\"),\n    @(\"2+80=\", \"76+14=\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $oldText\"\n    }\n}\n"}
